$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 6.704275
$ws.Range("H2").Value = 20.112825
$ws.Range("I2").Value = 0.4617710489234531
$ws.Range("J2").Value = 0.4617710489234532
$ws.Range("M2").Value = 0.1783253333333333
$ws.Range("N2").Value = 0.5349759999999999
$ws.Range("O2").Value = 0.03636495384128683
$ws.Range("P2").Value = 0.03636495384128683
$ws.Range("Q2").Value = 1.195542074133333
$ws.Range("R2").Value = 10.7598786672
$ws.Range("S2").Value = 0.01679228287934397
$ws.Range("T2").Value = 0.01679228287934398
$ws.Range("G3").Value = 6.704275
$ws.Range("H3").Value = 20.112825
$ws.Range("I3").Value = 0.4617710489234531
$ws.Range("J3").Value = 0.4617710489234532
$ws.Range("O3").Value = 0.7615778801978641
$ws.Range("P3").Value = 0.7615778801978641
$ws.Range("Q3").Value = 25.03779882355
$ws.Range("R3").Value = 225.34018941195
$ws.Range("S3").Value = 0.3516746165758676
$ws.Range("T3").Value = 0.3516746165758677
$ws.Range("G4").Value = 6.704275
$ws.Range("H4").Value = 20.112825
$ws.Range("I4").Value = 0.4617710489234531
$ws.Range("J4").Value = 0.4617710489234532
$ws.Range("O4").Value = 0.202057165960849
$ws.Range("P4").Value = 0.202057165960849
$ws.Range("Q4").Value = 6.642875014791667
$ws.Range("R4").Value = 59.785875133125
$ws.Range("S4").Value = 0.0933041494682415
$ws.Range("T4").Value = 0.0933041494682415
$ws.Range("I5").Value = 0.03922895479591048
$ws.Range("J5").Value = 0.03922895479591048
$ws.Range("M5").Value = 0.1783253333333333
$ws.Range("N5").Value = 0.5349759999999999
$ws.Range("O5").Value = 0.03636495384128683
$ws.Range("P5").Value = 0.03636495384128683
$ws.Range("Q5").Value = 0.1015651936
$ws.Range("R5").Value = 0.9140867423999999
$ws.Range("S5").Value = 0.001426559130395212
$ws.Range("T5").Value = 0.001426559130395212
$ws.Range("I6").Value = 0.03922895479591048
$ws.Range("J6").Value = 0.03922895479591048
$ws.Range("O6").Value = 0.7615778801978641
$ws.Range("P6").Value = 0.7615778801978641
$ws.Range("S6").Value = 0.02987590423584734
$ws.Range("T6").Value = 0.02987590423584734
$ws.Range("I7").Value = 0.03922895479591048
$ws.Range("J7").Value = 0.03922895479591048
$ws.Range("O7").Value = 0.202057165960849
$ws.Range("P7").Value = 0.202057165960849
$ws.Range("S7").Value = 0.007926491429667929
$ws.Range("T7").Value = 0.007926491429667929
$ws.Range("I8").Value = 0.4989999962806363
$ws.Range("J8").Value = 0.4989999962806364
$ws.Range("M8").Value = 0.1783253333333333
$ws.Range("N8").Value = 0.5349759999999999
$ws.Range("O8").Value = 0.03636495384128683
$ws.Range("P8").Value = 0.03636495384128683
$ws.Range("Q8").Value = 1.291929175587555
$ws.Range("R8").Value = 11.627362580288
$ws.Range("S8").Value = 0.01814611183154764
$ws.Range("T8").Value = 0.01814611183154764
$ws.Range("I9").Value = 0.4989999962806363
$ws.Range("J9").Value = 0.4989999962806364
$ws.Range("O9").Value = 0.7615778801978641
$ws.Range("P9").Value = 0.7615778801978641
$ws.Range("S9").Value = 0.3800273593861491
$ws.Range("T9").Value = 0.3800273593861491
$ws.Range("I10").Value = 0.4989999962806363
$ws.Range("J10").Value = 0.4989999962806364
$ws.Range("O10").Value = 0.202057165960849
$ws.Range("P10").Value = 0.202057165960849
$ws.Range("S10").Value = 0.1008265250629396
$ws.Range("T10").Value = 0.1008265250629396
